$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.773.57'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '3.761.39'
$ws.Range("E3").Value = '  -1.48%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '628.95'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").Value = '3.760.58'
$ws.Range("E7").Value = '  -1.44%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("E10").Value = '  -2.29%  '

$ws.Range("E11").Value = '  -0.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.93'
$ws.Range("E12").Value = '  +5.28%  '

$ws.Range("E13").Value = '  -4.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.75'
$ws.Range("E14").Value = '  -3.65%  '

$ws.Range("D15").Value = '4.393.82'
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").Value = '3.757.97'
$ws.Range("E16").Value = '  -1.84%  '

$ws.Range("D17").Value = '68.778.92'
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.60'
$ws.Range("E18").Value = '  -2.50%  '

$ws.Range("E20").Value = '  -2.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '466.84'
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.47'
$ws.Range("E22").Value = '  -2.18%  '

$ws.Range("E23").Value = '  -0.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.63'
$ws.Range("E24").Value = '  -2.69%  '

$ws.Range("E25").Value = '  -6.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.08'
$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.09'
$ws.Range("E28").Value = '  +0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").Value = '3.909.23'
$ws.Range("E30").Value = '  -1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.27'
$ws.Range("E31").Value = '  +1.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.65'
$ws.Range("E32").Value = '  -1.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.07'
$ws.Range("E33").Value = '  -2.92%  '

$ws.Range("E34").Value = '  +18.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.40'
$ws.Range("E35").Value = '  -2.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").Value = '3.713.84'
$ws.Range("E37").Value = '  -1.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.84'
$ws.Range("E38").Value = '  -2.66%  '

$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("E40").Value = '  -4.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.75'
$ws.Range("E41").Value = '  -2.68%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.958'
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("E45").Value = '  +4.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '155.59'
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.72'
$ws.Range("E47").Value = '  +3.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.84'
$ws.Range("E48").Value = '  +0.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.39'
$ws.Range("E49").Value = '  -3.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.292'
$ws.Range("E50").Value = '  -2.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.34'
$ws.Range("E51").Value = '  -1.29%  '
